# Daily attendance processing - 2025-12-17 10:32:25
#
# This script applies the daily attendance-processing edits to the
# "Session Analysis Results" sheet:
#   1. Swap the "Recorded By" order for every session that lists both
#      System and the instructor: "X, System" -> "System, X".
#   2. Recompute the top summary counters (Missing / Pending sessions).
#   3. Recompute the per-group statistics table (Missing / Pending cols).
#   4. Flip the 17/12/2025 GENERAL SURGERY sessions (and the matching
#      13/12/2025 SURGERY SEMINAR/SLIDE sessions) that were still
#      "Pending" from "Pending" to "Not Recorded" now that their date
#      has passed, re-styling those rows to the "Not Recorded" look.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. "<email>, System" -> "System, <email>" in the Recorded By column (G)
# ---------------------------------------------------------------------
$recordedByRows = @(2,3,4,23,24,25,26,45,46,47,48,67,68,69,70,89,90,91,110,111,112,131,132,133,152,153,154,173,174,175,194,195,196,197,216,217,218,219,238,239,240,241)

foreach ($r in $recordedByRows) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Value2 -eq "dnasr281@gmail.com, System") {
        $cell.Value = "System, dnasr281@gmail.com"
    }
}

# ---------------------------------------------------------------------
# 2. Top summary counters: Missing Sessions (L7) and Pending Sessions (L8)
# ---------------------------------------------------------------------
$ws.Range("L7").Value = 11
$ws.Range("L8").Value = 180

# ---------------------------------------------------------------------
# 3. Per-group statistics table: Missing (P) / Pending (Q) columns
# ---------------------------------------------------------------------
$groupStatRows = @(16,17,18,24,25,26)
foreach ($r in $groupStatRows) {
    $missing = $ws.Cells.Item($r, 16)
    $pending = $ws.Cells.Item($r, 17)
    $missing.Value = $missing.Value2 + 1
    $pending.Value = $pending.Value2 - 1
}

# ---------------------------------------------------------------------
# 4. Sessions that are now overdue flip from "Pending" to "Not Recorded"
# ---------------------------------------------------------------------
# Row 39 already displays "Not Recorded" with the correct red/pink
# styling - use it as the formatting template for every row that needs
# to look like a "Not Recorded" session.
$ws.Range("A39:I39").Copy() | Out-Null

# GENERAL SURGERY session 5 (17/12/2025) rows: text changes Pending -> Not Recorded
$generalSurgeryRows = @(27,49,71,198,220,242)
foreach ($r in $generalSurgeryRows) {
    $ws.Range("I$r").Value = "Not Recorded"
    $ws.Range("A" + $r + ":I" + $r).PasteSpecial(-4122) | Out-Null
}

# SURGERY SEMINAR/SLIDE rows that already say Not Recorded just get restyled
$seminarRows = @(61,210,232,254)
foreach ($r in $seminarRows) {
    $ws.Range("A" + $r + ":I" + $r).PasteSpecial(-4122) | Out-Null
}

$excel.CutCopyMode = 0
